$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.928.70'
$ws.Range("E2").Value = '  +0.23%  '

# Row 3
$ws.Range("D3").Value = '1.646.93'
$ws.Range("E3").Value = '  +1.02%  '

# Row 5
$ws.Range("D5").Value = '''216.47'
$ws.Range("E5").Value = '  +0.38%  '

# Row 6
$ws.Range("D6").Value = '''0.5064'
$ws.Range("E6").Value = '  -0.19%  '

# Row 7
$ws.Range("D7").Value = '''1.007'
$ws.Range("E7").Value = '  +0.57%  '

# Row 8
$ws.Range("D8").Value = '''0.2588'
$ws.Range("E8").Value = '  +0.51%  '

# Row 9
$ws.Range("D9").Value = '''0.06443'
$ws.Range("E9").Value = '  +1.74%  '

# Row 10
$ws.Range("D10").Value = '''20.49'
$ws.Range("E10").Value = '  +5.25%  '

# Row 11
$ws.Range("D11").Value = '''0.07820'
$ws.Range("E11").Value = '  +0.72%  '

# Row 12
$ws.Range("D12").Value = '''4.283'
$ws.Range("E12").Value = '  +0.73%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.646.45'
$ws.Range("E13").Value = '  +0.95%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.874.40'
$ws.Range("E14").Value = '  +1.03%  '

# Row 15
$ws.Range("D15").Value = '''0.5629'
$ws.Range("E15").Value = '  +2.20%  '

# Row 16
$ws.Range("D16").Value = '0.0₅7708'
$ws.Range("E16").Value = '  +0.76%  '

# Row 17
$ws.Range("D17").Value = '''63.47'
$ws.Range("E17").Value = '  -0.37%  '

# Row 18
$ws.Range("D18").Value = '25.965.25'
$ws.Range("E18").Value = '  +0.35%  '

# Row 19
$ws.Range("E19").Value = '  +0.52%  '

# Row 20
$ws.Range("D20").Value = '''193.56'
$ws.Range("E20").Value = '  -0.47%  '

# Row 21
$ws.Range("D21").Value = '''4.371'
$ws.Range("E21").Value = '  -0.75%  '

# Row 22
$ws.Range("D22").Value = '''9.956'
$ws.Range("E22").Value = '  +0.84%  '

# Row 23
$ws.Range("D23").Value = '''6.129'
$ws.Range("E23").Value = '  +1.77%  '

# Row 24
$ws.Range("D24").Value = '''1.007'
$ws.Range("E24").Value = '  +0.47%  '

# Row 25
$ws.Range("E25").Value = '  -5.73%  '

# Row 26
$ws.Range("D26").Value = '''141.53'
$ws.Range("E26").Value = '  -0.57%  '

# Row 27
$ws.Range("D27").Value = '''0.1240'
$ws.Range("E27").Value = '  -0.39%  '

# Row 28
$ws.Range("D28").Value = '''6.803'
$ws.Range("E28").Value = '  +0.66%  '

# Row 29
$ws.Range("D29").Value = '''15.56'
$ws.Range("E29").Value = '  -0.30%  '

# Row 30
$ws.Range("D30").Value = '''1.250'
$ws.Range("E30").Value = '  +0.97%  '

# Row 31
$ws.Range("D31").Value = '''0.04956'
$ws.Range("E31").Value = '  +1.48%  '

# Row 32
$ws.Range("D32").Value = '''3.304'
$ws.Range("E32").Value = '  +1.74%  '

# Row 33
$ws.Range("E33").Value = '  +1.72%  '

# Row 34
$ws.Range("D34").Value = '''1.576'
$ws.Range("E34").Value = '  +2.20%  '

# Row 35
$ws.Range("D35").Value = '''2.395'
$ws.Range("E35").Value = '  +1.17%  '

# Row 36
$ws.Range("D36").Value = '''0.9076'
$ws.Range("E36").Value = '  +1.32%  '

# Row 37
$ws.Range("D37").Value = '''0.5566'
$ws.Range("E37").Value = '  +0.82%  '

# Row 38
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.134.60'
$ws.Range("E38").Value = '  +1.62%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.556'
$ws.Range("E39").Value = '  +0.72%  '

# Row 40
$ws.Range("E40").Value = '  +1.03%  '

# Row 41
$ws.Range("D41").Value = '''1.006'
$ws.Range("E41").Value = '  +0.55%  '

# Row 42
$ws.Range("D42").Value = '''5.523'
$ws.Range("E42").Value = '  -1.50%  '

# Row 43
$ws.Range("D43").Value = '''0.8034'
$ws.Range("E43").Value = '  +0.68%  '

# Row 44
$ws.Range("D44").Value = '''98.90'
$ws.Range("E44").Value = '  +1.66%  '

# Row 45
$ws.Range("D45").Value = '1.784.91'
$ws.Range("E45").Value = '  +1.07%  '

# Row 46
$ws.Range("D46").Value = '0.0₈114'
$ws.Range("E46").Value = '  -3.96%  '

# Row 47
$ws.Range("D47").Value = '''55.79'
$ws.Range("E47").Value = '  +2.02%  '

# Row 48
$ws.Range("E48").Value = '  -3.32%  '

# Row 49
$ws.Range("D49").Value = '''7.763'
$ws.Range("E49").Value = '  +2.84%  '

# Row 50
$ws.Range("D50").Value = '''0.05046'
$ws.Range("E50").Value = '  -1.75%  '

# Row 51
$ws.Range("D51").Value = '''0.9995'
$ws.Range("E51").Value = '  -0.21%  '
